# Fix table header cells: "Tools" -> "Tool", "Example Paper(s)" -> "Representative study"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Tool"
$ws.Range("B1").Value = "Representative study"

$excel.Goto($ws.Range("A2"))
